$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Structural change: remove the "Makeup classes taken on January 9" column (O) ---
# Its counts are folded into the "January 2" column (N) instead, and the
# "Remaining Makeup classes" column (old P) shifts left to become the new O.
$ws.Columns("O:O").Delete()

# --- Fix up the view (diff removes topLeftCell and changes selection) ---
$ws.Range("I4").Select()

# --- Row 2: Jan-9 count (1) merged into Jan-2 count -> N2 becomes 2 ---
$ws.Range("N2").Value = 2

# --- Row 3: same merge ---
$ws.Range("N3").Value = 2

# --- Row 5: remove the "Late Notification" mark and the Jan-2 count ---
$ws.Range("H5").ClearContents()
$ws.Range("N5").ClearContents()

# --- Row 6: replace the "Late Notification" mark with a "CSE 208" mark ---
$ws.Range("H6").ClearContents()
$ws.Range("G6").Value = 1

# --- Row 7: remove the "CSE 208" mark and the Jan-2 count ---
$ws.Range("G7").ClearContents()
$ws.Range("N7").ClearContents()

# --- Recompute the "I" (total) column formulas as one shared formula ---
$ws.Range("I2:I7").Formula = "=C2+D2+E2+F2+G2+H2"

# --- Rebuild the "O" (remaining) column formulas (no longer subtracting the
#     removed January-9 column). O2 stands alone, O3:O7 share a formula,
#     matching the original sharing layout (which used to be P2 / P3:P7). ---
$ws.Range("O2").Formula = "=I2-J2-K2-L2-M2-N2"
$ws.Range("O3:O7").Formula = "=I3-J3-K3-L3-M3-N3"

$wb.Application.CalculateFull()
